# Update cryptos list: refresh prices, volume(1h) percentages, and fix
# two pairs of rows whose Coin/Link values were swapped relative to their
# correct price ordering (rows 38/39 and 43/44 and 47/48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.139.53'
$ws.Range("E2").Value = '  -4.11%  '
$ws.Range("D3").Value = '2.235.80'
$ws.Range("E3").Value = '  -4.75%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("E5").Value = '  -3.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  -6.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '68.36'
$ws.Range("E7").Value = '  -5.54%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.555'
$ws.Range("E9").Value = '  -5.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0972'
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.96'
$ws.Range("E11").Value = '  -0.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '34.78'
$ws.Range("E12").Value = '  +7.07%  '
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.65'
$ws.Range("E14").Value = '  -7.22%  '
$ws.Range("D15").Value = '2.565.43'
$ws.Range("E15").Value = '  -4.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.82'
$ws.Range("E16").Value = '  -8.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("E17").Value = '  -6.00%  '
$ws.Range("D18").Value = '2.241.19'
$ws.Range("E18").Value = '  -4.44%  '
$ws.Range("D19").Value = '41.900.24'
$ws.Range("E19").Value = '  -4.36%  '
$ws.Range("D20").Value = '0.0₃0965'
$ws.Range("E20").Value = '  -5.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.22'
$ws.Range("E21").Value = '  -6.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.86'
$ws.Range("E22").Value = '  -6.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.35'
$ws.Range("E23").Value = '  -7.77%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  -1.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.33'
$ws.Range("E27").Value = '  -6.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.94'
$ws.Range("E28").Value = '  -4.52%  '
$ws.Range("E29").Value = '  -3.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.59'
$ws.Range("E30").Value = '  -4.91%  '
$ws.Range("E31").Value = '  -8.61%  '
$ws.Range("E32").Value = '  -6.38%  '
$ws.Range("E33").Value = '  -7.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0709'
$ws.Range("E34").Value = '  -4.30%  '
$ws.Range("E35").Value = '  -2.82%  '
$ws.Range("E36").Value = '  -8.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.58'
$ws.Range("E37").Value = '  -4.39%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.24'
$ws.Range("E38").Value = '  -6.10%  '
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '20.97'
$ws.Range("E39").Value = '  +11.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.94'
$ws.Range("E40").Value = '  -7.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0262'
$ws.Range("E41").Value = '  -4.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.95'
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.80'
$ws.Range("E43").Value = '  -6.29%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.79'
$ws.Range("E44").Value = '  -4.84%  '
$ws.Range("E45").Value = '  -5.92%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("B47").Value = 'BitTorrent-New'
$ws.Range("C47").Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range("D47").Value = '0.0₃0157'
$ws.Range("E47").Value = '  +17.06%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.186'
$ws.Range("E48").Value = '  -6.15%  '
$ws.Range("E49").Value = '  -4.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  -6.61%  '
$ws.Range("E51").Value = '  +5.97%  '
